$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two header cells in row 1 (NO. PAJAK -> NO. FAKTUR, kode toko -> KODE TOKO)
$ws.Range("A9").Value = "Note:"
$ws.Range("A1").Value = "NO. FAKTUR"
$ws.Range("B1").Value = "KODE TOKO"
$ws.Range("A10").Value = "1. JIKA NO FAKTUR 24.0000024 MAKA HARUS DITAMBAHKAN SERI PAJAK DIDEPANNYA, SUPAYA SERI NO FAKTUR SAMA"

# Add a note row that should later be removed, highlighted in yellow
$ws.Range("A8").Value = "BAGIAN INI HAPUS "
$ws.Range("A8").Interior.Color = 65535

# Leave the selection on D8, matching the editor's final cursor position
$ws.Range("D8").Select()
